$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.405.41"
$ws.Range("E2").Value = "  -3.82%  "
$ws.Range("D3").Value = "1.571.87"
$ws.Range("E3").Value = "  -3.33%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("E5").Value = "  -0.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "289.55"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.71%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3674"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.26%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.31"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.44%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3401"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.171"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.03%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07645"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.59%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.27"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.074"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.31%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.922"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.82%  "
$ws.Range("D16").Value = "1.567.33"
$ws.Range("E16").Value = "  -3.63%  "
$ws.Range("E17").Value = "  -4.16%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "89.85"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.82%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06742"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.83%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.259"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.59"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.5314"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -7.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.25%  "
$ws.Range("D25").Value = "22.415.38"
$ws.Range("E25").Value = "  -3.84%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.354"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.911"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.99%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.04"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.70%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "146.34"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.976"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "125.79"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.24%  "
$ws.Range("D32").Value = "1.744.78"
$ws.Range("E32").Value = "  -3.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.253"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.019"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.017"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.60%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.14"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -9.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.08461"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.20%  "
$ws.Range("E38").Value = "  -3.90%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2325"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.49%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.529"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.92%  "
$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.06492"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.58%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.302"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.67%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.76"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.93%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6362"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.93%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.25"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.29%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9996"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5998"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.759"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.13%  "
$ws.Range("E49").Value = "  -5.15%  "
$ws.Range("E50").Value = "  +3.80%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "124.77"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.26%  "
